$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Restructure columns L:Q -> L:P ---
# Before: L=Test Accuracy, M=Model Size, N=Recall Drowsy, O=Recall Non-Drowsy, P=History Plot, Q=Confusion Matrix
# After:  L=Train Accuracy, M=Validation Accuracy, N=Model Size, O=History Plot, P=Confusion Matrix

# 1) Remove old "Model Size" column (M). This shifts N,O,P,Q left by one -> M,N,O,P
$ws.Columns("M:M").Delete()

# 2) Remove old "Recall Drowsy"/"Recall Non-Drowsy" columns (now at M:N). This shifts History Plot/Confusion Matrix left -> M,N
$ws.Columns("M:N").Delete()

# 3) Insert two fresh columns at M:N to host Validation Accuracy and Model Size
$ws.Columns("M:N").Insert()

# --- Header row ---
$ws.Range("L1").Value = "Train Accuracy"
$ws.Range("M1").Value = "Validation Accuracy"
$ws.Range("N1").Value = "Model Size"
$ws.Range("O1").Value = "History Plot"
$ws.Range("P1").Value = "Confusion Matrix"

# --- Data rows ---
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 1
$ws.Range("N2").Value = 189.1264686584473

$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.9666666388511658
$ws.Range("N3").Value = 189.1264686584473

$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 1
$ws.Range("N4").Value = 126.5555458068848

$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.9833333492279053
$ws.Range("N5").Value = 126.5555458068848

$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.9833333492279053
$ws.Range("N6").Value = 151.3449745178223

$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 1
$ws.Range("N7").Value = 151.3449745178223

$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 1
$ws.Range("N8").Value = 101.2740516662598

$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.9833333492279053
$ws.Range("N9").Value = 101.2740516662598
